# Insert a new weekly price record for "Perejil" (Vega Central Mapocho de
# Santiago) ahead of the existing row 384, pushing the remaining rows
# (384-419) down to (385-420) and growing the sheet's used range to R420.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 384..419 down to 385..420.
$ws.Rows.Item(384).Insert()

# Populate the newly inserted row 384 with the new observation.
$ws.Cells.Item(384, 1).Value  = 9
$ws.Cells.Item(384, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(384, 3).Value  = "Metropolitana"
$ws.Cells.Item(384, 4).Value  = 44783
$ws.Cells.Item(384, 5).Value  = 13
$ws.Cells.Item(384, 6).Value  = 100112044
$ws.Cells.Item(384, 7).Value  = "Perejil"
$ws.Cells.Item(384, 8).Value  = "Sin especificar"
$ws.Cells.Item(384, 9).Value  = "Primera"
$ws.Cells.Item(384, 10).Value = 70
$ws.Cells.Item(384, 11).Value = 15000
$ws.Cells.Item(384, 12).Value = 16000
$ws.Cells.Item(384, 13).Value = 15500
$ws.Cells.Item(384, 14).Value = "`$/docena de atados"
$ws.Cells.Item(384, 15).Value = "Región Metropolitana"
$ws.Cells.Item(384, 16).Value = 5167
$ws.Cells.Item(384, 17).Value = 3
$ws.Cells.Item(384, 18).Value = "Hortaliza"
